$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume(1h) columns store plain-text numeric- and percent-looking
# values (not real numbers) in the source sheet. Mark each target cell as
# Text *before* writing its new value so Excel keeps the literal string
# instead of auto-converting it into a formatted number/percentage.
$cells = @(
    "D2", "E2", "D3", "E3", "D4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9",
    "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15",
    "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D21", "E21", "D22", "E22",
    "D23", "E23", "D24", "E24", "D25", "E25", "E26", "E27", "D40", "E40", "D41", "E41",
    "D42", "E42", "E43", "D44", "E44", "D45", "E45", "E46", "D48", "E48", "D49", "E49",
    "D50", "E50"
)
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "261.23"
$ws.Range("E2").Value = "0.00%"
$ws.Range("D3").Value = "26.79"
$ws.Range("E3").Value = "-1.91%"
$ws.Range("D4").Value = "4.723"
$ws.Range("D5").Value = "0.06209"
$ws.Range("E5").Value = "1.97%"
$ws.Range("D6").Value = "6.724"
$ws.Range("E6").Value = "0.78%"
$ws.Range("D7").Value = "0.8496"
$ws.Range("E7").Value = "0.42%"
$ws.Range("D8").Value = "0.9130"
$ws.Range("E8").Value = "-1.35%"
$ws.Range("D9").Value = "0.1407"
$ws.Range("E9").Value = "0.01%"
$ws.Range("D10").Value = "0.04976"
$ws.Range("E10").Value = "5.43%"
$ws.Range("D11").Value = "0.07093"
$ws.Range("E11").Value = "-0.19%"
$ws.Range("D12").Value = "0.03116"
$ws.Range("E12").Value = "1.04%"
$ws.Range("D13").Value = "0.09048"
$ws.Range("E13").Value = "-0.13%"
$ws.Range("D14").Value = "0.001539"
$ws.Range("E14").Value = "0.34%"
$ws.Range("D15").Value = "0.0006161"
$ws.Range("E15").Value = "1.33%"
$ws.Range("D16").Value = "0.005952"
$ws.Range("E16").Value = "-3.08%"
$ws.Range("D17").Value = "3.449"
$ws.Range("E17").Value = "-0.01%"
$ws.Range("D18").Value = "3.172"
$ws.Range("E18").Value = "1.02%"
$ws.Range("D19").Value = "2.168"
$ws.Range("E19").Value = "0.22%"
$ws.Range("D21").Value = "0.1311"
$ws.Range("E21").Value = "1.74%"
$ws.Range("D22").Value = "4.094"
$ws.Range("E22").Value = "0.26%"
$ws.Range("D23").Value = "0.04236"
$ws.Range("E23").Value = "-0.23%"
$ws.Range("D24").Value = "0.001180"
$ws.Range("E24").Value = "-3.46%"
$ws.Range("D25").Value = "0.004063"
$ws.Range("E25").Value = "3.88%"
$ws.Range("E26").Value = "0.02%"
$ws.Range("E27").Value = "4.11%"
$ws.Range("D40").Value = "0.03936"
$ws.Range("E40").Value = "1.57%"
$ws.Range("D41").Value = "0.1112"
$ws.Range("E41").Value = "-0.24%"
$ws.Range("D42").Value = "0.004120"
$ws.Range("E42").Value = "0.41%"
$ws.Range("E43").Value = "-3.34%"
$ws.Range("D44").Value = "0.01327"
$ws.Range("E44").Value = "-18.50%"
$ws.Range("D45").Value = "0.00005167"
$ws.Range("E45").Value = "0.43%"
$ws.Range("E46").Value = "0.11%"
$ws.Range("D48").Value = "0.2483"
$ws.Range("E48").Value = "82.79%"
$ws.Range("D49").Value = "0.00002103"
$ws.Range("E49").Value = "0.11%"
$ws.Range("D50").Value = "0.0002003"
$ws.Range("E50").Value = "0.11%"
